$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (A3:T3) values each increment by 1 (21-40 -> 22-41)
for ($col = 1; $col -le 20; $col++) {
    $cell = $ws.Cells.Item(3, $col)
    $cell.Value2 = $cell.Value2 + 1
}

# Remove the extra styled rows 19, 20, 21 (each previously had A19/A20/A21 with style "1")
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()

# Keep row 4 registered as part of the worksheet's used range (so the sheet
# dimension stays A1:T4 instead of collapsing once the trailing rows are gone)
$ws.Cells.Item(4, 1).Style = "Normal"

# Update the selection to match the new used range
$ws.Range("A4:T20").Select()

$wb.Save()
